$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "Log" sheet (Process.VaR/cVaR rows) ------------------------
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

# Row 14: mark the original VaR/cVaR script entry COMPLETED and extend the
# comment in column B to mention the stress-testing follow-up work.
$log.Cells.Item(14, 4).Value = "COMPLETED"
$log.Cells.Item(14, 2).Value = "Make a script that measures VaR and cVaR for an equity so I know how it works. Made another one with stress testing."

# Row 15: new note about scrolling through CSVs for portfolio VaR/cVaR + GARCH.
$log.Cells.Item(15, 2).Value = 'Make another file that scrolls through all .CSVs to calculate "Portfolio VaR and cVaR" a (e.g. GARCH)nd understand improved mathematical concepts'
$log.Rows.Item(15).RowHeight = 29

# ---------------------------------------------------------------------------
# 2) Add the new "PromptsToFollowUpOn" sheet at the end of the workbook ----
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PromptsToFollowUpOn"

# A1 holds the literal Warren-Buffet prompt text; a leading apostrophe marks
# it as a text literal (quote-prefixed) the way a pasted/typed prompt would
# be, and wrap text keeps the long paragraph readable.
$newSheet.Cells.Item(1, 1).Value = "'Now pretend you are Warren Buffet. I want to invest in Chinese stocks but want to apply Warren's mindset of long term purchases to a foreign market. Help me understand how Warren would evaluate equities and what considerations he would take into account before buying and his timing of purchases."
$newSheet.Cells.Item(1, 1).WrapText = $true

$newSheet.Cells.Item(1, 2).Value = "Provided the answer I expected"

$newSheet.Range("B2").Select()

# Row 16: new note about the portfolio-wide / individual stats + GARCH script.
# (Written after the new sheet's prompts so the shared-string table keeps the
# same allocation order as the authored edit.)
$log.Cells.Item(16, 2).Value = "Created another one that scrolls through all CSVs in folder to calculate portfolio wide and invdidual stats and stress tests. Created a script that forecasts GARCH but doesn't return anything useful. Might have to break it down into pieces."
$log.Rows.Item(16).RowHeight = 29

# Move the visible selection to where the user left off.
$log.Range("B17").Select()

# Keep "Log" as the active/visible tab, matching the original workbook.
$log.Activate()
